$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add 5 new per-sample toggle columns ("run_silva", "run_greengenes",
# "run_rdp", "run_crest", "run_ncbi_blast") right before the existing
# "taxa_barstacks" / "max_taxa" columns, keeping the sheet's usual blank
# separator column between field groups.
# ------------------------------------------------------------------

# Insert 6 blank columns (CD:CI). The existing "taxa_barstacks"/"max_taxa"
# columns (old CD/CE) get pushed right to CJ/CK; CI is left as the blank
# separator column (matching the existing CC/BY/... pattern elsewhere on
# this sheet).
$ws.Range("CD1:CI7").Insert(-4161)

# Copy the look (font/fill/alignment/number format) of each header/body
# row from the column that used to be CD (now shifted to CJ) onto the
# new columns, so they match the sheet's existing styling for that row
# instead of picking up a blank default style.
$ws.Range("CJ1").Copy()
$ws.Range("CD1:CH1").PasteSpecial(-4122)

$ws.Range("CJ2").Copy()
$ws.Range("CD2:CH2").PasteSpecial(-4122)

$ws.Range("CJ3").Copy()
$ws.Range("CD3:CH7").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# --- Column widths for the 5 new columns -----------------------------
# COM ColumnWidth is offset from the raw OOXML column width by a fixed
# 5/6 padding; compensate so the saved width matches the target.
$pad = 5/6
$ws.Range("CD1").ColumnWidth = 18.5 - $pad
$ws.Range("CE1").ColumnWidth = 19 - $pad
$ws.Range("CF1").ColumnWidth = 18.5 - $pad
$ws.Range("CG1").ColumnWidth = 18.83203125 - $pad
$ws.Range("CH1").ColumnWidth = 20.33203125 - $pad

# --- Row 1 (long-form question headers) -------------------------------
$ws.Range("CD1").Value = "When this option is turned on, we will run the Silva taxonomic classification and produce a report. This option can be set to ""TRUE"" or to ""FALSE""."
$ws.Range("CE1").Value = "When this option is turned on, we will run the Green Genes taxonomic classification and produce a report. This option can be set to ""TRUE"" or to ""FALSE""."
$ws.Range("CF1").Value = "When this option is turned on, we will run the RDP taxonomic classification and produce a report. This option can be set to ""TRUE"" or to ""FALSE""."
$ws.Range("CG1").Value = "When this option is turned on, we will run the CREST LCA taxonomic classification and produce a report. This option can be set to ""TRUE"" or to ""FALSE""."
$ws.Range("CH1").Value = "When this option is turned on, we will run the NCBI BLAST taxonomic classification and produce a report. This option can be set to ""TRUE"" or to ""FALSE""."

# --- Row 2 (short column/field names) ----------------------------------
$ws.Range("CD2").Value = "run_silva"
$ws.Range("CE2").Value = "run_greengenes"
$ws.Range("CF2").Value = "run_rdp"
$ws.Range("CG2").Value = "run_crest"
$ws.Range("CH2").Value = "run_ncbi_blast"

# --- Rows 3-7 (per-sample boolean values) -------------------------------
for ($r = 3; $r -le 7; $r++) {
    $ws.Range("CD$r").Value = $true
    $ws.Range("CE$r").Value = $true
    $ws.Range("CF$r").Value = $true
    $ws.Range("CG$r").Value = $false
    $ws.Range("CH$r").Value = $false
}

# --- Mirror the author's final cursor position after adding the columns
[void]$ws.Range("CF11").Select()

Write-Output "done"
